$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 27 and 28 swap coin content (ranking positions exchanged)
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1331"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -10.43%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'18.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.17%  "

# Per-row updated Price (D) and Volume 1h (E) values
$ws.Range("D2").Value = "30.242.49"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "1.927.43"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'246.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").Value = "'0.7171"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -11.46%  "
$ws.Range("D7").Value = "'0.9997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.3257"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.28%  "
$ws.Range("D9").Value = "'26.46"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "'0.06819"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.8027"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").Value = "'0.07940"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.928.49"
$ws.Range("E13").Value = "  -2.98%  "
$ws.Range("D14").Value = "'5.399"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'94.38"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.27%  "
$ws.Range("E16").Value = "  +3.89%  "
$ws.Range("D17").Value = "'260.66"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "30.243.55"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").Value = "'5.843"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'0.000007939"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "2.180.96"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'0.9999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'6.863"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("D25").Value = "'9.668"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "'160.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D29").Value = "'2.285"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("D30").Value = "'1.361"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").Value = "'4.196"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "'0.05072"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'1.197"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "'0.7417"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("D37").Value = "'2.726"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").Value = "'2.809"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("D40").Value = "'79.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.27%  "
$ws.Range("D41").Value = "'6.563"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "'0.4451"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("D43").Value = "'2.002"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").Value = "'0.8331"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "'102.57"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'9.729"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").Value = "'36.22"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'1.487"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").Value = "'0.4106"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.55%  "
